$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 184.6
$ws.Range("I33").Value = 170.35
$ws.Range("K33").Value = 170.35
$ws.Range("M33").Value = 58.65000000000001
$ws.Range("H112").Value = 2383.9
$ws.Range("I112").Value = 750
$ws.Range("J112").Value = 2565.4443
$ws.Range("K112").Value = 2250
$ws.Range("L112").Value = 7696.3329
$ws.Range("M112").Value = -1142
$ws.Range("N112").Value = -9912.332900000001
$ws.Range("H116").Value = 3152.6667
$ws.Range("I116").Value = 3088
$ws.Range("J116").Value = 3379
$ws.Range("K116").Value = 3088
$ws.Range("L116").Value = 3379
$ws.Range("M116").Value = 354
$ws.Range("N116").Value = -10263
$ws.Range("H138").Value = 547520.4399999999
$ws.Range("I138").Value = 4937.421
$ws.Range("J138").Value = 929338.1
$ws.Range("K138").Value = 14812.263
$ws.Range("L138").Value = 2788014.3
$ws.Range("M138").Value = -9672.263000000001
$ws.Range("N138").Value = -2798294.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1309.8125
$ws.Range("I2").Value = 996.1818
$ws.Range("K2").Value = 996.1818
$ws.Range("M2").Value = -883.1818
$ws.Range("H32").Value = 4393.8604
$ws.Range("I32").Value = 3644.8171
$ws.Range("J32").Value = 19749.25
$ws.Range("K32").Value = 3644.8171
$ws.Range("L32").Value = 19749.25
$ws.Range("M32").Value = -3357.8171
$ws.Range("N32").Value = -20323.25
$ws.Range("H45").Value = 29663.312
$ws.Range("I45").Value = 86592.8
$ws.Range("K45").Value = 86592.8
$ws.Range("M45").Value = -86215.8
$ws.Range("H63").Value = 5465.4165
$ws.Range("I63").Value = 3635.625
$ws.Range("J63").Value = 9125
$ws.Range("K63").Value = 3635.625
$ws.Range("L63").Value = 9125
$ws.Range("M63").Value = -2949.625
$ws.Range("N63").Value = -10497
$ws.Range("H66").Value = 5465.4165
$ws.Range("I66").Value = 3635.625
$ws.Range("J66").Value = 9125
$ws.Range("K66").Value = 18178.125
$ws.Range("L66").Value = 45625
$ws.Range("M66").Value = -14746.125
$ws.Range("N66").Value = -52489
$ws.Range("H74").Value = 157357.56
$ws.Range("I74").Value = 244648.12
$ws.Range("J74").Value = 2920.3845
$ws.Range("K74").Value = 244648.12
$ws.Range("L74").Value = 2920.3845
$ws.Range("M74").Value = -243774.12
$ws.Range("N74").Value = -4668.3845
$ws.Range("H77").Value = 157357.56
$ws.Range("I77").Value = 244648.12
$ws.Range("J77").Value = 2920.3845
$ws.Range("K77").Value = 1223240.6
$ws.Range("L77").Value = 14601.9225
$ws.Range("M77").Value = -1218872.6
$ws.Range("N77").Value = -23337.9225
$ws.Range("H116").Value = 1309.8125
$ws.Range("I116").Value = 996.1818
$ws.Range("K116").Value = 996.1818
$ws.Range("M116").Value = 1297.8182
$ws.Range("H122").Value = 5128.1904
$ws.Range("I122").Value = 4120.125
$ws.Range("J122").Value = 8354
$ws.Range("K122").Value = 12360.375
$ws.Range("L122").Value = 25062
$ws.Range("M122").Value = -9910.375
$ws.Range("N122").Value = -29962

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1309.8125
$ws.Range("I3").Value = 996.1818
$ws.Range("K3").Value = 996.1818
$ws.Range("M3").Value = -882.1818
$ws.Range("H94").Value = 41667468
$ws.Range("I94").Value = 55556104
$ws.Range("K94").Value = 55556104
$ws.Range("M94").Value = -55555653
$ws.Range("H107").Value = 1709.8928
$ws.Range("I107").Value = 1457.8125
$ws.Range("J107").Value = 2046
$ws.Range("K107").Value = 1457.8125
$ws.Range("L107").Value = 2046
$ws.Range("M107").Value = 462.1875
$ws.Range("N107").Value = -5886
$ws.Range("H134").Value = 2176.125
$ws.Range("I134").Value = 1543.7307
$ws.Range("K134").Value = 4631.1921
$ws.Range("M134").Value = -2096.1921

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3461.2837
$ws.Range("I31").Value = 2737.2964
$ws.Range("K31").Value = 2737.2964
$ws.Range("M31").Value = -2442.2964
$ws.Range("H34").Value = 3461.2837
$ws.Range("I34").Value = 2737.2964
$ws.Range("K34").Value = 2737.2964
$ws.Range("M34").Value = -2535.2964

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 2216.5
$ws.Range("I16").Value = 399
$ws.Range("J16").Value = 2580
$ws.Range("K16").Value = 1197
$ws.Range("L16").Value = 7740
$ws.Range("M16").Value = -1024
$ws.Range("N16").Value = -8086
$ws.Range("H92").Value = 475.7143
$ws.Range("I92").Value = 491.5
$ws.Range("J92").Value = 381
$ws.Range("K92").Value = 1474.5
$ws.Range("L92").Value = 1143
$ws.Range("M92").Value = -226.5
$ws.Range("N92").Value = -3639
$ws.Range("H137").Value = 4357.4707
$ws.Range("I137").Value = 3937.0833
$ws.Range("J137").Value = 5366.4
$ws.Range("K137").Value = 11811.2499
$ws.Range("L137").Value = 16099.2
$ws.Range("M137").Value = -6711.249899999999
$ws.Range("N137").Value = -26299.2
$ws.Range("H140").Value = 10324.088
$ws.Range("I140").Value = 5455.4546
$ws.Range("K140").Value = 16366.3638
$ws.Range("M140").Value = -11186.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 90912430
$ws.Range("I80").Value = 166669390
$ws.Range("K80").Value = 166669390
$ws.Range("M80").Value = -166668392
$ws.Range("H83").Value = 90912430
$ws.Range("I83").Value = 166669390
$ws.Range("K83").Value = 833346950
$ws.Range("M83").Value = -833341958
$ws.Range("H97").Value = 452.27274
$ws.Range("I97").Value = 478.375
$ws.Range("K97").Value = 478.375
$ws.Range("M97").Value = 17.625
$ws.Range("H113").Value = 4388.517
$ws.Range("I113").Value = 4213.1113
$ws.Range("J113").Value = 6756.5
$ws.Range("K113").Value = 4213.1113
$ws.Range("L113").Value = 6756.5
$ws.Range("M113").Value = -2043.1113
$ws.Range("N113").Value = -11096.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2241.8
$ws.Range("I61").Value = 2241.8
$ws.Range("K61").Value = 2241.8
$ws.Range("M61").Value = -2039.8
$ws.Range("H113").Value = 2241.8
$ws.Range("I113").Value = 2241.8
$ws.Range("K113").Value = 2241.8
$ws.Range("M113").Value = -71.80000000000018
$ws.Range("H136").Value = 6374.4
$ws.Range("I136").Value = 5401.857
$ws.Range("J136").Value = 19990
$ws.Range("K136").Value = 16205.571
$ws.Range("L136").Value = 59970
$ws.Range("M136").Value = -13655.571
$ws.Range("N136").Value = -65070

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 913.4286
$ws.Range("I113").Value = 1065.2727
$ws.Range("J113").Value = 746.4
$ws.Range("K113").Value = 3195.8181
$ws.Range("L113").Value = 2239.2
$ws.Range("M113").Value = -1025.8181
$ws.Range("N113").Value = -6579.2
$ws.Range("H122").Value = 20837176
$ws.Range("I122").Value = 3565.5557
$ws.Range("J122").Value = 83338000
$ws.Range("K122").Value = 10696.6671
$ws.Range("L122").Value = 250014000
$ws.Range("M122").Value = -8246.667099999999
$ws.Range("N122").Value = -250018900
